$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.24752933333333
$ws.Range("H2").Value = 150.742588
$ws.Range("I2").Value = 0.1294604271951564
$ws.Range("J2").Value = 0.1294604271951564
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.701354
$ws.Range("N2").Value = 8.104061999999999
$ws.Range("O2").Value = 0.02221077311549548
$ws.Range("P2").Value = 0.02221077311549548
$ws.Range("Q2").Value = 135.7363643547173
$ws.Range("R2").Value = 1221.627279192456
$ws.Range("S2").Value = 0.002875416175866739
$ws.Range("T2").Value = 0.002875416175866739

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.24752933333333
$ws.Range("H3").Value = 150.742588
$ws.Range("I3").Value = 0.1294604271951564
$ws.Range("J3").Value = 0.1294604271951564
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("N3").Value = 264.449356
$ws.Range("O3").Value = 0.7247753838328104
$ws.Range("P3").Value = 0.7247753838328105
$ws.Range("Q3").Value = 4429.308924263702
$ws.Range("R3").Value = 39863.78031837332
$ws.Range("S3").Value = 0.09382973081152905
$ws.Range("T3").Value = 0.09382973081152908

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.24752933333333
$ws.Range("H4").Value = 150.742588
$ws.Range("I4").Value = 0.1294604271951564
$ws.Range("J4").Value = 0.1294604271951564
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.24063
$ws.Range("N4").Value = 0.72189
$ws.Range("O4").Value = 0.001978481285600361
$ws.Range("P4").Value = 0.001978481285600361
$ws.Range("Q4").Value = 12.09106298348
$ws.Range("R4").Value = 108.81956685132
$ws.Range("S4").Value = 0.0002561350324314449
$ws.Range("T4").Value = 0.000256135032431445

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.24752933333333
$ws.Range("H5").Value = 150.742588
$ws.Range("I5").Value = 0.1294604271951564
$ws.Range("J5").Value = 0.1294604271951564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.53182233333333
$ws.Range("N5").Value = 91.595467
$ws.Range("O5").Value = 0.2510353617660938
$ws.Range("P5").Value = 0.2510353617660938
$ws.Range("Q5").Value = 1534.148638294288
$ws.Range("R5").Value = 13807.33774464859
$ws.Range("S5").Value = 0.03249914517532913
$ws.Range("T5").Value = 0.03249914517532913

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 151.42276
$ws.Range("H6").Value = 454.26828
$ws.Range("I6").Value = 0.3901337131747328
$ws.Range("J6").Value = 0.3901337131747328
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.701354
$ws.Range("N6").Value = 8.104061999999999
$ws.Range("O6").Value = 0.02221077311549548
$ws.Range("P6").Value = 0.02221077311549548
$ws.Range("Q6").Value = 409.04647841704
$ws.Range("R6").Value = 3681.41830575336
$ws.Range("S6").Value = 0.008665171388029779
$ws.Range("T6").Value = 0.008665171388029779

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 151.42276
$ws.Range("H7").Value = 454.26828
$ws.Range("I7").Value = 0.3901337131747328
$ws.Range("J7").Value = 0.3901337131747328
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("N7").Value = 264.449356
$ws.Range("O7").Value = 0.7247753838328104
$ws.Range("P7").Value = 0.7247753838328105
$ws.Range("Q7").Value = 13347.88378858085
$ws.Range("R7").Value = 120130.9540972277
$ws.Range("S7").Value = 0.2827593117123365
$ws.Range("T7").Value = 0.2827593117123365

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 151.42276
$ws.Range("H8").Value = 454.26828
$ws.Range("I8").Value = 0.3901337131747328
$ws.Range("J8").Value = 0.3901337131747328
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.24063
$ws.Range("N8").Value = 0.72189
$ws.Range("O8").Value = 0.001978481285600361
$ws.Range("P8").Value = 0.001978481285600361
$ws.Range("Q8").Value = 36.43685873880001
$ws.Range("R8").Value = 327.9317286492
$ws.Range("S8").Value = 0.0007718722503979879
$ws.Range("T8").Value = 0.0007718722503979879

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 151.42276
$ws.Range("H9").Value = 454.26828
$ws.Range("I9").Value = 0.3901337131747328
$ws.Range("J9").Value = 0.3901337131747328
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.53182233333333
$ws.Range("N9").Value = 91.595467
$ws.Range("O9").Value = 0.2510353617660938
$ws.Range("P9").Value = 0.2510353617660938
$ws.Range("Q9").Value = 4623.212805542974
$ws.Range("R9").Value = 41608.91524988676
$ws.Range("S9").Value = 0.09793735782396852
$ws.Range("T9").Value = 0.09793735782396852

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 85.147481
$ws.Range("H10").Value = 255.442443
$ws.Range("I10").Value = 0.2193785328573129
$ws.Range("J10").Value = 0.2193785328573129
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.701354
$ws.Range("N10").Value = 8.104061999999999
$ws.Range("O10").Value = 0.02221077311549548
$ws.Range("P10").Value = 0.02221077311549548
$ws.Range("Q10").Value = 230.013488389274
$ws.Range("R10").Value = 2070.121395503466
$ws.Range("S10").Value = 0.004872566819704048
$ws.Range("T10").Value = 0.004872566819704048

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 85.147481
$ws.Range("H11").Value = 255.442443
$ws.Range("I11").Value = 0.2193785328573129
$ws.Range("J11").Value = 0.2193785328573129
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("N11").Value = 264.449356
$ws.Range("O11").Value = 0.7247753838328104
$ws.Range("P11").Value = 0.7247753838328105
$ws.Range("Q11").Value = 7505.732171824078
$ws.Range("R11").Value = 67551.5895464167
$ws.Range("S11").Value = 0.1590001603563378
$ws.Range("T11").Value = 0.1590001603563378

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 85.147481
$ws.Range("H12").Value = 255.442443
$ws.Range("I12").Value = 0.2193785328573129
$ws.Range("J12").Value = 0.2193785328573129
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.24063
$ws.Range("N12").Value = 0.72189
$ws.Range("O12").Value = 0.001978481285600361
$ws.Range("P12").Value = 0.001978481285600361
$ws.Range("Q12").Value = 20.48903835303
$ws.Range("R12").Value = 184.40134517727
$ws.Range("S12").Value = 0.0004340363217206576
$ws.Range("T12").Value = 0.0004340363217206576

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 85.147481
$ws.Range("H13").Value = 255.442443
$ws.Range("I13").Value = 0.2193785328573129
$ws.Range("J13").Value = 0.2193785328573129
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.53182233333333
$ws.Range("N13").Value = 91.595467
$ws.Range("O13").Value = 0.2510353617660938
$ws.Range("P13").Value = 0.2510353617660938
$ws.Range("Q13").Value = 2599.707762022876
$ws.Range("R13").Value = 23397.36985820588
$ws.Range("S13").Value = 0.05507176935955044
$ws.Range("T13").Value = 0.05507176935955044

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 101.3126446666667
$ws.Range("H14").Value = 303.937934
$ws.Range("I14").Value = 0.2610273267727979
$ws.Range("J14").Value = 0.2610273267727979
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.701354
$ws.Range("N14").Value = 8.104061999999999
$ws.Range("O14").Value = 0.02221077311549548
$ws.Range("P14").Value = 0.02221077311549548
$ws.Range("Q14").Value = 273.6813179208787
$ws.Range("R14").Value = 2463.131861287908
$ws.Range("S14").Value = 0.005797618731894915
$ws.Range("T14").Value = 0.005797618731894915

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 101.3126446666667
$ws.Range("H15").Value = 303.937934
$ws.Range("I15").Value = 0.2610273267727979
$ws.Range("J15").Value = 0.2610273267727979
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("N15").Value = 264.449356
$ws.Range("O15").Value = 0.7247753838328104
$ws.Range("P15").Value = 0.7247753838328105
$ws.Range("Q15").Value = 8930.687878918945
$ws.Range("R15").Value = 80376.19091027051
$ws.Range("S15").Value = 0.1891861809526071
$ws.Range("T15").Value = 0.1891861809526071

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 101.3126446666667
$ws.Range("H16").Value = 303.937934
$ws.Range("I16").Value = 0.2610273267727979
$ws.Range("J16").Value = 0.2610273267727979
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.24063
$ws.Range("N16").Value = 0.72189
$ws.Range("O16").Value = 0.001978481285600361
$ws.Range("P16").Value = 0.001978481285600361
$ws.Range("Q16").Value = 24.37886168614001
$ws.Range("R16").Value = 219.40975517526
$ws.Range("S16").Value = 0.0005164376810502709
$ws.Range("T16").Value = 0.0005164376810502709

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 101.3126446666667
$ws.Range("H17").Value = 303.937934
$ws.Range("I17").Value = 0.2610273267727979
$ws.Range("J17").Value = 0.2610273267727979
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.53182233333333
$ws.Range("N17").Value = 91.595467
$ws.Range("O17").Value = 0.2510353617660938
$ws.Range("P17").Value = 0.2510353617660938
$ws.Range("Q17").Value = 3093.259667082798
$ws.Range("R17").Value = 27839.33700374518
$ws.Range("S17").Value = 0.06552708940724572
$ws.Range("T17").Value = 0.06552708940724572
